# This workbook is a daily-updated COVID-19 tracking sheet for Valais.
# Columns B, H, J and K are driven by shared formulas (cumulative totals /
# simple sums) and recompute automatically; only the raw input columns
# (C, E, F, G, L, M) need to be edited here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 232 - revise hospitalised-outside-ICU figure
$ws.Range("G232").Value = 91

# Row 244 - extra-hospital deaths revised
$ws.Range("M244").Value = 4

# Row 255 - hospital / extra-hospital death split revised
$ws.Range("L255").Value = 6
$ws.Range("M255").Value = 7

# Row 256 - new positive cases (hors SI) revised
$ws.Range("F256").Value = 21

# Row 257 - new positive cases (hors SI) revised
$ws.Range("F257").Value = 24

# Row 258 - new cases, hors SI, hospital deaths revised
$ws.Range("C258").Value = 554
$ws.Range("F258").Value = 26
$ws.Range("L258").Value = 8

# Row 259 - intubated patients revised
$ws.Range("E259").Value = 36

# Row 260 - intubated patients / extra-hospital deaths revised
$ws.Range("E260").Value = 35
$ws.Range("M260").Value = 10

# Row 261 - new cases / intubated patients revised
$ws.Range("C261").Value = 319
$ws.Range("E261").Value = 37

# Row 262 - new cases / intubated patients revised
$ws.Range("C262").Value = 303
$ws.Range("E262").Value = 34

# Row 263 - intubated patients revised
$ws.Range("E263").Value = 33

# Row 264 - new cases / intubated patients revised
$ws.Range("C264").Value = 119
$ws.Range("E264").Value = 34

# Row 265 - intubated patients revised
$ws.Range("E265").Value = 33

# Row 266 - new cases / intubated patients / hospital deaths revised
$ws.Range("C266").Value = 280
$ws.Range("E266").Value = 34
$ws.Range("L266").Value = 7

# Row 267 - new cases / intubated patients / extra-hospital deaths revised
$ws.Range("C267").Value = 214
$ws.Range("E267").Value = 34
$ws.Range("M267").Value = 9

# Row 268 - new cases / intubated patients / hospitalised hors SI / hospital deaths revised
$ws.Range("C268").Value = 141
$ws.Range("E268").Value = 37
$ws.Range("G268").Value = 221
$ws.Range("L268").Value = 2

# Row 269 - newly filled-in day (previously blank inputs)
$ws.Range("C269").Value = 13
$ws.Range("E269").Value = 35
$ws.Range("F269").Value = 26
$ws.Range("G269").Value = 212
$ws.Range("L269").Value = 0
$ws.Range("M269").Value = 0

$wb.Application.CalculateFullRebuild()
